$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Voltage Table")

# Swap the stepping-table values between row 19 ("Outer Rejection Electrode")
# and row 20 ("Inner Rejection Electrode") so that row 19 now holds -3500
# and row 20 now holds 4000. Row 35/36 recompute automatically since they
# are driven by formulas referencing row 19/20.
$ws.Range("H19").Value = -3500
$ws.Range("H20").Value = 4000

# Update the active selection to match the saved view state.
$ws.Range("H21").Select()
